$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")

# F2: availability for HH2_BIO_SR_C_NEW (input HH2_CU) in 2020
$ws.Range("F2").Value = 34.05888

# Row 3: HH2_NGA_CL_CCS_NEW (input HH2_CT)
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 810.3494229786608
$ws.Range("L3").Value = 1501.629429327508

# Row 4: HH2_COA_CL_CCS_NEW (input HH2_CT)
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0.7152364799999911
